$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Network Scanner" section -------------------------------------------
$ws.Range("A35").Value = "Network Scanner"
$ws.Range("B36").Value = "NetDiscover"
$ws.Range("C37").Value = "netdiscover -r XXX.XXX.X.1/24"

# --- "Intro to ARP" section -----------------------------------------------
$ws.Range("B39").Value = "Intro to ARP"
$ws.Range("C41").Value = "Sends ARP Requests."
$ws.Range("C42").Value = "Receives ARP Responses."
$ws.Range("C40").Value = "ARP is a protocol. A device can broadcast ARP Requests for MAC addresses by IP, and receive a resopnse back."

# --- "Create Network Scanner with Scapy" section --------------------------
$ws.Range("B43").Value = "Create Network Scanner with Scapy"
$ws.Range("C44").Value = "scapy.arping(ip)"

# --- "Crate ARP Packet" section -------------------------------------------
$ws.Range("C46").Value = "Create ARP Packet."
$ws.Range("B45").Value = "Crate ARP Packet"
$ws.Range("C47").Value = "request = scapy.ARP()"
$ws.Range("C48").Value = "broadcast = scapy.Ether()"
$ws.Range("C49").Value = "arp_packet = broadcast/request"

# Update the current selection to reflect where the author left off editing.
$ws.Range("C50").Select()
